$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the second calibration table (Coolant Temp Sensor) row: the
# 20C measurement moved to 21C and its resistance reading changed to 6400.
$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 6400

# Leave the sheet's selection where the author last left it.
$ws.Range("C16").Select()
